$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.851.74'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.236.83'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.88%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.09'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.38'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -5.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.572'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -3.10%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -7.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.72'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -8.48%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.51%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -7.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.575.72'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.85%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -4.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.231.40'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.00'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -4.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.688.99'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.07'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -8.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0966'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.67'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.12'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.99'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -7.29%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -8.29%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.10'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.94%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.55'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -6.60%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -9.45%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.93'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -3.03%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '156.93'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0832'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -6.05%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.29'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.18%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.65'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.60%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -9.36%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.92%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '15.77'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.08%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -9.45%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -11.23%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -6.36%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.703.25'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '82.51'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.27%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -6.79%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.74%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '101.51'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '71.44'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -5.20%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '56.22'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -6.24%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.61'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.04%  '
